# Auto-generated Excel COM-interop script
# Applies numeric cell updates to the Maduin_Profits profit-calculator sheets
# (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR) per the scheduled runner refresh.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 1642.3334
$ws.Range("I9").Value = 1099.5
$ws.Range("J9").Value = 1797.4286
$ws.Range("K9").Value = 1099.5
$ws.Range("L9").Value = 1797.4286
$ws.Range("M9").Value = -930.5
$ws.Range("N9").Value = -2135.4286

$ws.Range("H32").Value = 1435.875
$ws.Range("I32").Value = 1495
$ws.Range("J32").Value = 1416.1666
$ws.Range("K32").Value = 1495
$ws.Range("L32").Value = 1416.1666
$ws.Range("M32").Value = -1169
$ws.Range("N32").Value = -2068.1666

$ws.Range("H40").Value = 3333.3333
$ws.Range("I40").Value = 0
$ws.Range("J40").Value = 3333.3333
$ws.Range("K40").Value = 0
$ws.Range("L40").Value = 3333.3333
$ws.Range("M40").ClearContents()
$ws.Range("N40").Value = -3683.3333

$ws.Range("H55").Value = 314.33334
$ws.Range("I55").Value = 315
$ws.Range("J55").Value = 314
$ws.Range("K55").Value = 315
$ws.Range("L55").Value = 314
$ws.Range("M55").Value = -101
$ws.Range("N55").Value = -742

$ws.Range("H112").Value = 1923.3889
$ws.Range("I112").Value = 1299.3334
$ws.Range("J112").Value = 2048.2
$ws.Range("K112").Value = 3898.0002
$ws.Range("L112").Value = 6144.599999999999
$ws.Range("M112").Value = -2790.0002
$ws.Range("N112").Value = -8360.599999999999

$ws.Range("H114").Value = 69832.5
$ws.Range("I114").Value = 0
$ws.Range("J114").Value = 69832.5
$ws.Range("K114").Value = 0
$ws.Range("L114").Value = 69832.5
$ws.Range("N114").Value = -78510.5

$ws.Range("H132").Value = 3427.1177
$ws.Range("I132").Value = 2329.8096
$ws.Range("J132").Value = 5199.6924
$ws.Range("K132").Value = 6989.4288
$ws.Range("L132").Value = 15599.0772
$ws.Range("M132").Value = -4459.4288
$ws.Range("N132").Value = -20659.0772

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3495.7222
$ws.Range("I32").Value = 2581.818
$ws.Range("J32").Value = 13548.667
$ws.Range("K32").Value = 2581.818
$ws.Range("L32").Value = 13548.667
$ws.Range("M32").Value = -2294.818
$ws.Range("N32").Value = -14122.667

$ws.Range("H122").Value = 3259.4614
$ws.Range("I122").Value = 3246
$ws.Range("J122").Value = 3304.3333
$ws.Range("K122").Value = 9738
$ws.Range("L122").Value = 9912.999899999999
$ws.Range("M122").Value = -7288
$ws.Range("N122").Value = -14812.9999

$ws.Range("H132").Value = 1044
$ws.Range("I132").Value = 1201.5
$ws.Range("J132").Value = 99
$ws.Range("K132").Value = 3604.5
$ws.Range("L132").Value = 297
$ws.Range("M132").Value = -1074.5
$ws.Range("N132").Value = -5357

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 6060.857
$ws.Range("I86").Value = 4582
$ws.Range("J86").Value = 11483.333
$ws.Range("K86").Value = 4582
$ws.Range("L86").Value = 11483.333
$ws.Range("M86").Value = -3459
$ws.Range("N86").Value = -13729.333

$ws.Range("H89").Value = 6060.857
$ws.Range("I89").Value = 4582
$ws.Range("J89").Value = 11483.333
$ws.Range("K89").Value = 22910
$ws.Range("L89").Value = 57416.665
$ws.Range("M89").Value = -17294
$ws.Range("N89").Value = -68648.66500000001

$ws.Range("H134").Value = 2474.6155
$ws.Range("I134").Value = 2297.2727
$ws.Range("J134").Value = 3450
$ws.Range("K134").Value = 6891.8181
$ws.Range("L134").Value = 10350
$ws.Range("M134").Value = -4356.8181
$ws.Range("N134").Value = -15420

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 933.76666
$ws.Range("I58").Value = 727.9231
$ws.Range("J58").Value = 2271.75
$ws.Range("K58").Value = 727.9231
$ws.Range("L58").Value = 2271.75
$ws.Range("M58").Value = -524.9231
$ws.Range("N58").Value = -2677.75

$ws.Range("H74").Value = 24999
$ws.Range("I74").Value = 0
$ws.Range("J74").Value = 24999
$ws.Range("K74").Value = 0
$ws.Range("L74").Value = 24999
$ws.Range("N74").Value = -26747

$ws.Range("H77").Value = 24999
$ws.Range("I77").Value = 0
$ws.Range("J77").Value = 24999
$ws.Range("K77").Value = 0
$ws.Range("L77").Value = 74997
$ws.Range("N77").Value = -83733

$ws.Range("H122").Value = 1852.5883
$ws.Range("I122").Value = 1122.875
$ws.Range("J122").Value = 2501.2222
$ws.Range("K122").Value = 3368.625
$ws.Range("L122").Value = 7503.6666
$ws.Range("M122").Value = -918.625
$ws.Range("N122").Value = -12403.6666

$ws.Range("H132").Value = 2698.5908
$ws.Range("I132").Value = 3015.111
$ws.Range("J132").Value = 1274.25
$ws.Range("K132").Value = 9045.332999999999
$ws.Range("L132").Value = 3822.75
$ws.Range("M132").Value = -6515.332999999999
$ws.Range("N132").Value = -8882.75

$ws.Range("H136").Value = 933.76666
$ws.Range("I136").Value = 727.9231
$ws.Range("J136").Value = 2271.75
$ws.Range("K136").Value = 2183.7693
$ws.Range("L136").Value = 6815.25
$ws.Range("M136").Value = 366.2307000000001
$ws.Range("N136").Value = -11915.25

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 784.5714
$ws.Range("I107").Value = 529
$ws.Range("J107").Value = 854.2727
$ws.Range("K107").Value = 1587
$ws.Range("L107").Value = 2562.8181
$ws.Range("M107").Value = 333
$ws.Range("N107").Value = -6402.8181

$ws.Range("H140").Value = 2285.6
$ws.Range("I140").Value = 708.8889
$ws.Range("J140").Value = 4650.6665
$ws.Range("K140").Value = 2126.6667
$ws.Range("L140").Value = 13951.9995
$ws.Range("M140").Value = 3053.3333
$ws.Range("N140").Value = -24311.9995

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 1054.4166
$ws.Range("I122").Value = 1085.4
$ws.Range("J122").Value = 899.5
$ws.Range("K122").Value = 3256.2
$ws.Range("L122").Value = 2698.5
$ws.Range("M122").Value = -806.2000000000003
$ws.Range("N122").Value = -7598.5

$ws.Range("H132").Value = 2553.5386
$ws.Range("I132").Value = 2098.5
$ws.Range("J132").Value = 8014
$ws.Range("K132").Value = 6295.5
$ws.Range("L132").Value = 24042
$ws.Range("M132").Value = -3765.5
$ws.Range("N132").Value = -29102

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1096
$ws.Range("I22").Value = 200
$ws.Range("J22").Value = 1992
$ws.Range("K22").Value = 200
$ws.Range("L22").Value = 1992
$ws.Range("M22").Value = 95
$ws.Range("N22").Value = -2582

$ws.Range("H27").Value = 1096
$ws.Range("I27").Value = 200
$ws.Range("J27").Value = 1992
$ws.Range("K27").Value = 200
$ws.Range("L27").Value = 1992
$ws.Range("M27").Value = -93
$ws.Range("N27").Value = -2206

$ws.Range("H38").Value = 40000
$ws.Range("I38").Value = 0
$ws.Range("J38").Value = 40000
$ws.Range("K38").Value = 0
$ws.Range("L38").Value = 40000
$ws.Range("N38").Value = -40820

$ws.Range("H55").Value = 441.6154
$ws.Range("I55").Value = 90.25
$ws.Range("J55").Value = 597.7778
$ws.Range("K55").Value = 90.25
$ws.Range("L55").Value = 597.7778
$ws.Range("M55").Value = 82.75
$ws.Range("N55").Value = -943.7778

$ws.Range("H76").Value = 16500
$ws.Range("I76").Value = 0
$ws.Range("J76").Value = 16500
$ws.Range("K76").Value = 0
$ws.Range("L76").Value = 16500
$ws.Range("N76").Value = -17176

$ws.Range("H79").Value = 16500
$ws.Range("I79").Value = 0
$ws.Range("J79").Value = 16500
$ws.Range("K79").Value = 0
$ws.Range("L79").Value = 16500
$ws.Range("N79").Value = -18840

$ws.Range("H122").Value = 6446.625
$ws.Range("I122").Value = 5872.5
$ws.Range("J122").Value = 7020.75
$ws.Range("K122").Value = 17617.5
$ws.Range("L122").Value = 21062.25
$ws.Range("M122").Value = -15167.5
$ws.Range("N122").Value = -25962.25

$ws.Range("H132").Value = 4347.647
$ws.Range("I132").Value = 4593.933
$ws.Range("J132").Value = 2500.5
$ws.Range("K132").Value = 13781.799
$ws.Range("L132").Value = 7501.5
$ws.Range("M132").Value = -11251.799
$ws.Range("N132").Value = -12561.5

$ws.Range("H136").Value = 3599.3333
$ws.Range("I136").Value = 2913.4285
$ws.Range("J136").Value = 6000
$ws.Range("K136").Value = 8740.2855
$ws.Range("L136").Value = 18000
$ws.Range("M136").Value = -6190.2855
$ws.Range("N136").Value = -23100

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H74").Value = 35000
$ws.Range("I74").Value = 35000
$ws.Range("J74").Value = 0
$ws.Range("K74").Value = 35000
$ws.Range("L74").Value = 0
$ws.Range("M74").Value = -34064

$ws.Range("H77").Value = 35000
$ws.Range("I77").Value = 35000
$ws.Range("J77").Value = 0
$ws.Range("K77").Value = 105000
$ws.Range("L77").Value = 0
$ws.Range("M77").Value = -100320

$ws.Range("H122").Value = 975
$ws.Range("I122").Value = 975
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 2925
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -475

$ws.Range("H126").Value = 3042.6924
$ws.Range("I126").Value = 3042.6924
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 9128.0772
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -6658.0772
$ws.Range("N126").ClearContents()

$ws.Range("H132").Value = 1230.5625
$ws.Range("I132").Value = 1336.6364
$ws.Range("J132").Value = 997.2
$ws.Range("K132").Value = 4009.9092
$ws.Range("L132").Value = 2991.6
$ws.Range("M132").Value = -1479.9092
$ws.Range("N132").Value = -8051.6

$ws.Range("H136").Value = 1002
$ws.Range("I136").Value = 1002
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 3006
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -456
